$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Body (word/document.xml): remove the first three paragraphs
#    ("A simple demonstration of a query :", the field-code paragraph
#    showing the invalid expression error, and "End of demonstration.")
#    leaving only the bookmark paragraph that was already present.
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$thirdPara = $d.Paragraphs.Item(3)
$deleteRange = $d.Range($firstPara.Range.Start, $thirdPara.Range.End)
$deleteRange.Delete()

# ------------------------------------------------------------------
# 2. Header (word/header1.xml): merge the run-split, spell-checked
#    "A simple demonstration of a query :" text into two runs with an
#    explicit en-US language and a grammar-check proofErr pair instead
#    of the two spell-check proofErr pairs.
# ------------------------------------------------------------------
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrRange = $hdr.Range
$hdrRange.SetRange(0, 36)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">A simple demonstration of a </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>query :</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '</w:p>'

$hdrRange.InsertXML($newParaXml)
